$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B9 held "Klein et al. (2018)" (bold). Move that value+bold formatting up to
# B4 (which previously held the now-removed "In press" placeholder), then
# clear B9 entirely so the cell goes away.
$citation = $ws.Range("B9").Value2

$ws.Range("B4").Value = $citation
$ws.Range("B4").Font.Bold = $true

$ws.Range("B9").Clear()

# Update the active selection to reflect where editing left off.
$ws.Range("D10").Select()
